$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry data to be rearranged between rows: D, K, L, M, N, O, P, Q, R, S, T
$cols = @("D","K","L","M","N","O","P","Q","R","S","T")

# Capture current (before) values for every affected row so we can freely
# reassign them without clobbering data we still need to read.
$rowsToRead = @(2,4,5,6,7,8,10,11,12,13)
$orig = @{}
foreach ($r in $rowsToRead) {
    $vals = @()
    foreach ($c in $cols) {
        $vals += ,($ws.Range("$c$r").Value2)
    }
    $orig[$r] = $vals
}

# Mapping: target row -> source row (data that should end up in target row)
$mapping = @{
    5  = 2
    2  = 4
    8  = 5
    4  = 6
    6  = 7
    11 = 8
    12 = 10
    13 = 11
    10 = 12
    7  = 13
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $vals = $orig[$sourceRow]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$targetRow").Value2 = $vals[$i]
    }
}
